# ---------------------------------------------------------------------------
# Update GraficiRisultati/30_3/20MB/results.xlsx:
#   - rename the two bandwidth-series headers to include units "(Mb/s)"
#   - add two new columns with the total duration (s) of each download/upload
#     test, with header text + three rows of data
#   - widen the columns to fit the new/renamed headers
#   - update the chart title (profile + file-size subtitle) and the two
#     series names so the chart reflects the new header text
#   - move/resize the chart so it sits below the data table instead of next
#     to it
#   - leave the active selection on F4, matching the author's last edit
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename existing headers to include the unit of measure -------------
$ws.Range("B1").Value = "Banda in download (Mb/s)"
$ws.Range("C1").Value = "Banda in upload (Mb/s)"

# --- 2. Add the two new "total time" columns --------------------------------
$ws.Range("D1").Value = "Tempo totale download (s)"
$ws.Range("E1").Value = "Tempo totale upload (s)"

$ws.Range("D2").Value = 52.003
$ws.Range("E2").Value = 22.037

$ws.Range("D3").Value = 19.002
$ws.Range("E3").Value = 19.036

$ws.Range("D4").Value = 19.002
$ws.Range("E4").Value = 16.047

# --- 3. Widen the columns so the longer headers are readable ----------------
$ws.Columns.Item(1).ColumnWidth = 22.498697916666668   # -> stored width 23.33203125
$ws.Columns.Item(2).ColumnWidth = 22.330729166666668   # -> stored width 23.1640625
$ws.Columns.Item(3).ColumnWidth = 22.330729166666668   # -> stored width 23.1640625
$ws.Columns.Item(4).ColumnWidth = 23.166666666666668   # -> stored width 24
$ws.Columns.Item(5).ColumnWidth = 22.330729166666668   # -> stored width 23.1640625

# --- 4. Update the chart: title, series names -------------------------------
$co = $ws.ChartObjects(1)
$chart = $co.Chart

$nl = [char]10
$chart.ChartTitle.Text = "Profilo: 30/3" + $nl + "Dimensione file: 20MB"

$s1 = $chart.SeriesCollection(1)
$s2 = $chart.SeriesCollection(2)
$s1.Name = "Banda in download (Mb/s)"
$s2.Name = "Banda in upload (Mb/s)"

# --- 5. Reposition the chart below the table (was to the right) -------------
# Target anchor (EMU, twoCellAnchor): from col0/off12700,row8/off0
#                                     to   col11/off622300,row29/off63500
# With the new column widths above (126.25/125.375/125.375/129.75/125.375 pt
# for A:E, 58.4375pt default beyond) and the default 16pt row height, that
# anchor corresponds to this Left/Top/Width/Height (in points):
$co.Left = 1
$co.Top = 128
$co.Width = 1030.75
$co.Height = 341

# --- 6. Restore the author's final selection --------------------------------
[void]$ws.Range("F4").Select()
